$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1617.5454
$ws.Range("I6").Value = 448
$ws.Range("J6").Value = 2056.125
$ws.Range("K6").Value = 1344
$ws.Range("L6").Value = 6168.375
$ws.Range("M6").Value = -1232
$ws.Range("N6").Value = -6392.375
$ws.Range("H32").Value = 1939.1111
$ws.Range("J32").Value = 1764.5714
$ws.Range("L32").Value = 1764.5714
$ws.Range("N32").Value = -2416.5714
$ws.Range("H53").Value = 12324.667
$ws.Range("I53").Value = 27599.25
$ws.Range("K53").Value = 27599.25
$ws.Range("M53").Value = -26962.25
$ws.Range("H98").Value = 6660.375
$ws.Range("I98").Value = 6680.5
$ws.Range("J98").Value = 6600
$ws.Range("K98").Value = 6680.5
$ws.Range("L98").Value = 6600
$ws.Range("M98").Value = -5182.5
$ws.Range("N98").Value = -9596
$ws.Range("H122").Value = 6660.375
$ws.Range("I122").Value = 6680.5
$ws.Range("J122").Value = 6600
$ws.Range("K122").Value = 20041.5
$ws.Range("L122").Value = 19800
$ws.Range("M122").Value = -17591.5
$ws.Range("N122").Value = -24700
$ws.Range("H129").Value = 1133.6316
$ws.Range("I129").Value = 497
$ws.Range("J129").Value = 1150.8379
$ws.Range("K129").Value = 1491
$ws.Range("L129").Value = 3452.5137
$ws.Range("M129").Value = 3509
$ws.Range("N129").Value = -13452.5137
$ws.Range("H132").Value = 863.60315
$ws.Range("I132").Value = 805.2069
$ws.Range("J132").Value = 1541
$ws.Range("K132").Value = 2415.6207
$ws.Range("L132").Value = 4623
$ws.Range("M132").Value = 114.3793000000001
$ws.Range("N132").Value = -9683
$ws.Range("H137").Value = 51407.35
$ws.Range("I137").Value = 1117.75
$ws.Range("J137").Value = 84933.75
$ws.Range("K137").Value = 3353.25
$ws.Range("L137").Value = 254801.25
$ws.Range("M137").Value = -803.25
$ws.Range("N137").Value = -259901.25
$ws.Range("H138").Value = 1689.899
$ws.Range("I138").Value = 1015.91174
$ws.Range("J138").Value = 2042.4462
$ws.Range("K138").Value = 3047.73522
$ws.Range("L138").Value = 6127.3386
$ws.Range("M138").Value = 2092.26478
$ws.Range("N138").Value = -16407.3386

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4604.879
$ws.Range("I32").Value = 2884.6667
$ws.Range("K32").Value = 2884.6667
$ws.Range("M32").Value = -2597.6667
$ws.Range("H45").Value = 1337.7916
$ws.Range("I45").Value = 1154.6666
$ws.Range("J45").Value = 1643
$ws.Range("K45").Value = 1154.6666
$ws.Range("L45").Value = 1643
$ws.Range("M45").Value = -777.6666
$ws.Range("N45").Value = -2397
$ws.Range("H61").Value = 23969
$ws.Range("I61").Value = 25952.625
$ws.Range("K61").Value = 25952.625
$ws.Range("M61").Value = -25740.625
$ws.Range("H74").Value = 721.6177
$ws.Range("I74").Value = 567.7273
$ws.Range("K74").Value = 567.7273
$ws.Range("M74").Value = 306.2727
$ws.Range("H77").Value = 721.6177
$ws.Range("I77").Value = 567.7273
$ws.Range("K77").Value = 2838.6365
$ws.Range("M77").Value = 1529.3635
$ws.Range("H102").Value = 1399.4
$ws.Range("I102").Value = 1399.4
$ws.Range("K102").Value = 1399.4
$ws.Range("M102").Value = 222.5999999999999
$ws.Range("H132").Value = 1412.494
$ws.Range("I132").Value = 1073.9056
$ws.Range("J132").Value = 2010.6666
$ws.Range("K132").Value = 3221.7168
$ws.Range("L132").Value = 6031.9998
$ws.Range("M132").Value = -691.7168000000001
$ws.Range("N132").Value = -11091.9998
$ws.Range("H136").Value = 23969
$ws.Range("I136").Value = 25952.625
$ws.Range("K136").Value = 77857.875
$ws.Range("M136").Value = -75307.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812
$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808
$ws.Range("H94").Value = 889.625
$ws.Range("I94").Value = 615.6
$ws.Range("K94").Value = 615.6
$ws.Range("M94").Value = -164.6
$ws.Range("H134").Value = 3687.3333
$ws.Range("I134").Value = 3837.4167
$ws.Range("K134").Value = 11512.2501
$ws.Range("M134").Value = -8977.250100000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1972.7894
$ws.Range("I31").Value = 1320.091
$ws.Range("K31").Value = 1320.091
$ws.Range("M31").Value = -1025.091
$ws.Range("H34").Value = 1972.7894
$ws.Range("I34").Value = 1320.091
$ws.Range("K34").Value = 1320.091
$ws.Range("M34").Value = -1118.091
$ws.Range("H58").Value = 1450719.5
$ws.Range("I58").Value = 2718924.2
$ws.Range("J58").Value = 1342.7142
$ws.Range("K58").Value = 2718924.2
$ws.Range("L58").Value = 1342.7142
$ws.Range("M58").Value = -2718721.2
$ws.Range("N58").Value = -1748.7142
$ws.Range("H99").Value = 3142.9
$ws.Range("J99").Value = 3241.5
$ws.Range("L99").Value = 3241.5
$ws.Range("N99").Value = -6237.5
$ws.Range("H126").Value = 3142.9
$ws.Range("J126").Value = 3241.5
$ws.Range("L126").Value = 9724.5
$ws.Range("N126").Value = -14664.5
$ws.Range("H132").Value = 1445.7241
$ws.Range("I132").Value = 1085.75
$ws.Range("J132").Value = 2245.6667
$ws.Range("K132").Value = 3257.25
$ws.Range("L132").Value = 6737.000100000001
$ws.Range("M132").Value = -727.25
$ws.Range("N132").Value = -11797.0001
$ws.Range("H136").Value = 1450719.5
$ws.Range("I136").Value = 2718924.2
$ws.Range("J136").Value = 1342.7142
$ws.Range("K136").Value = 8156772.600000001
$ws.Range("L136").Value = 4028.1426
$ws.Range("M136").Value = -8154222.600000001
$ws.Range("N136").Value = -9128.142599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 250.44444
$ws.Range("J2").Value = 221.8
$ws.Range("L2").Value = 1330.8
$ws.Range("N2").Value = -1556.8
$ws.Range("H107").Value = 450.35294
$ws.Range("I107").Value = 364
$ws.Range("K107").Value = 1092
$ws.Range("M107").Value = 828
$ws.Range("H113").Value = 55975.35
$ws.Range("I113").Value = 183967.33
$ws.Range("J113").Value = 1121.6428
$ws.Range("K113").Value = 551901.99
$ws.Range("L113").Value = 3364.9284
$ws.Range("M113").Value = -549731.99
$ws.Range("N113").Value = -7704.928400000001
$ws.Range("H140").Value = 3753.0588
$ws.Range("I140").Value = 2835.3333
$ws.Range("K140").Value = 8505.999899999999
$ws.Range("M140").Value = -3325.999899999999
$ws.Range("H141").Value = 2968.9412
$ws.Range("I141").Value = 2629.2666
$ws.Range("K141").Value = 7887.7998
$ws.Range("M141").Value = -2707.7998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3973.5
$ws.Range("J80").Value = 3973.5
$ws.Range("L80").Value = 3973.5
$ws.Range("N80").Value = -5969.5
$ws.Range("H83").Value = 3973.5
$ws.Range("J83").Value = 3973.5
$ws.Range("L83").Value = 19867.5
$ws.Range("N83").Value = -29851.5
$ws.Range("H97").Value = 1427.579
$ws.Range("I97").Value = 1506.75
$ws.Range("J97").Value = 1291.8572
$ws.Range("K97").Value = 1506.75
$ws.Range("L97").Value = 1291.8572
$ws.Range("M97").Value = -1010.75
$ws.Range("N97").Value = -2283.8572
$ws.Range("H113").Value = 1246.1111
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 1403
$ws.Range("K113").Value = 1050
$ws.Range("L113").Value = 1403
$ws.Range("M113").Value = 1120
$ws.Range("N113").Value = -5743
$ws.Range("H126").Value = 2461806.5
$ws.Range("I126").Value = 2780832.5
$ws.Range("K126").Value = 8342497.5
$ws.Range("M126").Value = -8340027.5
$ws.Range("H132").Value = 1167412.2
$ws.Range("I132").Value = 1480834.2
$ws.Range("J132").Value = 3273.2856
$ws.Range("K132").Value = 4442502.6
$ws.Range("L132").Value = 9819.856800000001
$ws.Range("M132").Value = -4439972.6
$ws.Range("N132").Value = -14879.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2352.7778
$ws.Range("I46").Value = 1366.6666
$ws.Range("J46").Value = 2845.8333
$ws.Range("K46").Value = 1366.6666
$ws.Range("L46").Value = 2845.8333
$ws.Range("M46").Value = -1178.6666
$ws.Range("N46").Value = -3221.8333
$ws.Range("H55").Value = 550.7143
$ws.Range("I55").Value = 451
$ws.Range("J55").Value = 683.6667
$ws.Range("K55").Value = 451
$ws.Range("L55").Value = 683.6667
$ws.Range("M55").Value = -278
$ws.Range("N55").Value = -1029.6667
$ws.Range("H132").Value = 3126.7673
$ws.Range("I132").Value = 2555.074
$ws.Range("J132").Value = 4091.5
$ws.Range("K132").Value = 7665.222
$ws.Range("L132").Value = 12274.5
$ws.Range("M132").Value = -5135.222
$ws.Range("N132").Value = -17334.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -31020
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696
$ws.Range("H122").Value = 35331.418
$ws.Range("I122").Value = 42002.25
$ws.Range("J122").Value = 1977.25
$ws.Range("K122").Value = 126006.75
$ws.Range("L122").Value = 5931.75
$ws.Range("M122").Value = -123556.75
$ws.Range("N122").Value = -10831.75
$ws.Range("H132").Value = 15434.983
$ws.Range("I132").Value = 1024.3721
$ws.Range("J132").Value = 48048.473
$ws.Range("K132").Value = 3073.1163
$ws.Range("L132").Value = 144145.419
$ws.Range("M132").Value = -543.1163000000001
$ws.Range("N132").Value = -149205.419
$ws.Range("H136").Value = 15874296
$ws.Range("I136").Value = 24155532
$ws.Range("K136").Value = 72466596
$ws.Range("M136").Value = -72464046
